$d = $word.ActiveDocument

# --- First paragraph formatting ---
$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right) reserving 5pt of space,
# without turning on a visible line style/color (matches <w:pBdr><w:top w:space="5"/>...).
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Left indent: 120 twips (6pt) -> 225 twips (11.25pt)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# --- Remove the trailing " " run at the end of paragraph 1 ---
$r1 = $p1.Range
$spaceRange = $d.Range($r1.End - 2, $r1.End - 1)
$spaceRange.Delete()

# --- Update the remaining run's text ---
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$searchRange = $d.Range($r1.Start, $r1.End)
$searchRange.Find.Execute("**ID__AFFARS_pgi_5349_topic_2__ID**", $false, $false, $false, $false, $false, `
                           $true, 1, $false, "**ID__AFFARS_AF_PGI_5349__ID**", 2)
